$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H2").Value = 1365.3636
$ws.Range("I2").Value = 1131.2858
$ws.Range("J2").Value = 1775
$ws.Range("K2").Value = 1131.2858
$ws.Range("L2").Value = 1775
$ws.Range("M2").Value = -1018.2858
$ws.Range("N2").Value = -2001

$ws.Range("H5").Value = 39
$ws.Range("I5").Value = 39
$ws.Range("K5").Value = 39
$ws.Range("M5").Value = 76

$ws.Range("H53").Value = 352.5
$ws.Range("I53").Value = 494.76923
$ws.Range("J53").Value = 184.36363
$ws.Range("K53").Value = 494.76923
$ws.Range("L53").Value = 184.36363
$ws.Range("M53").Value = 142.23077
$ws.Range("N53").Value = -1458.36363

$ws.Range("H111").Value = 400
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 400
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 1200
$ws.Range("M111").ClearContents()
$ws.Range("N111").Value = -7334

$ws.Range("H138").Value = 3783.9583
$ws.Range("J138").Value = 3792.1738
$ws.Range("L138").Value = 11376.5214
$ws.Range("N138").Value = -21656.5214

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H16").Value = 18729.285
$ws.Range("I16").Value = 25276.25
$ws.Range("J16").Value = 10000
$ws.Range("K16").Value = 25276.25
$ws.Range("L16").Value = 10000
$ws.Range("M16").Value = -24989.25
$ws.Range("N16").Value = -10574

$ws.Range("H74").Value = 1635.6111
$ws.Range("I74").Value = 1143.7646
$ws.Range("K74").Value = 1143.7646
$ws.Range("M74").Value = -269.7646

$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

$ws.Range("H77").Value = 1635.6111
$ws.Range("I77").Value = 1143.7646
$ws.Range("K77").Value = 5718.823
$ws.Range("M77").Value = -1350.823

$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

$ws.Range("H88").Value = 8986.375
$ws.Range("J88").Value = 9555.857
$ws.Range("L88").Value = 9555.857
$ws.Range("N88").Value = -10367.857

$ws.Range("H91").Value = 8986.375
$ws.Range("J91").Value = 9555.857
$ws.Range("L91").Value = 9555.857
$ws.Range("N91").Value = -12363.857

$ws.Range("H132").Value = 10004
$ws.Range("I132").Value = 10004
$ws.Range("K132").Value = 30012
$ws.Range("M132").Value = -27482

$ws = $wb.Sheets.Item("BSM")
$ws.Range("H86").Value = 2007
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()

$ws.Range("H89").Value = 2007
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()

$ws.Range("H105").Value = 3435.923
$ws.Range("I105").Value = 2381
$ws.Range("K105").Value = 2381
$ws.Range("M105").Value = -634

$ws.Range("H107").Value = 1831.1111
$ws.Range("I107").Value = 1211.4286
$ws.Range("K107").Value = 1211.4286
$ws.Range("M107").Value = 708.5714

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H31").Value = 3050.5264
$ws.Range("I31").Value = 2282.4546
$ws.Range("J31").Value = 4106.625
$ws.Range("K31").Value = 2282.4546
$ws.Range("L31").Value = 4106.625
$ws.Range("M31").Value = -1987.4546
$ws.Range("N31").Value = -4696.625

$ws.Range("H34").Value = 3050.5264
$ws.Range("I34").Value = 2282.4546
$ws.Range("J34").Value = 4106.625
$ws.Range("K34").Value = 2282.4546
$ws.Range("L34").Value = 4106.625
$ws.Range("M34").Value = -2080.4546
$ws.Range("N34").Value = -4510.625

$ws.Range("H58").Value = 4998.5
$ws.Range("I58").Value = 4998.5
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 4998.5
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -4795.5
$ws.Range("N58").ClearContents()

$ws.Range("H120").Value = 19999.5
$ws.Range("I120").Value = 19999
$ws.Range("K120").Value = 19999
$ws.Range("M120").Value = -16370

$ws.Range("H132").Value = 1643.5454
$ws.Range("I132").Value = 1342.2222
$ws.Range("K132").Value = 4026.6666
$ws.Range("M132").Value = -1496.6666

$ws.Range("H134").Value = 3974.1875
$ws.Range("I134").Value = 4042.2144
$ws.Range("J134").Value = 3498
$ws.Range("K134").Value = 12126.6432
$ws.Range("L134").Value = 10494
$ws.Range("M134").Value = -9591.643199999999
$ws.Range("N134").Value = -15564

$ws.Range("H136").Value = 4998.5
$ws.Range("I136").Value = 4998.5
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 14995.5
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -12445.5
$ws.Range("N136").ClearContents()

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H2").Value = 153.5
$ws.Range("I2").Value = 83.09090999999999
$ws.Range("J2").Value = 264.14285
$ws.Range("K2").Value = 498.5454599999999
$ws.Range("L2").Value = 1584.8571
$ws.Range("M2").Value = -385.5454599999999
$ws.Range("N2").Value = -1810.8571

$ws.Range("H23").Value = 531.63635
$ws.Range("J23").Value = 504.9
$ws.Range("L23").Value = 1514.7
$ws.Range("N23").Value = -1984.7

$ws.Range("H37").Value = 75000
$ws.Range("J37").Value = 75000
$ws.Range("L37").Value = 225000
$ws.Range("N37").Value = -225224

$ws.Range("H50").Value = 431.5
$ws.Range("I50").Value = 552.3333
$ws.Range("J50").Value = 69
$ws.Range("K50").Value = 1656.9999
$ws.Range("L50").Value = 207
$ws.Range("M50").Value = -1175.9999
$ws.Range("N50").Value = -1169

$ws.Range("H53").Value = 431.5
$ws.Range("I53").Value = 552.3333
$ws.Range("J53").Value = 69
$ws.Range("K53").Value = 1656.9999
$ws.Range("L53").Value = 207
$ws.Range("M53").Value = -1175.9999
$ws.Range("N53").Value = -1169

$ws.Range("H114").Value = 257.8
$ws.Range("I114").Value = 329.66666
$ws.Range("J114").Value = 150
$ws.Range("K114").Value = 988.9999799999999
$ws.Range("L114").Value = 450
$ws.Range("M114").Value = 2265.00002
$ws.Range("N114").Value = -6958

$ws.Range("H117").Value = 712.4286
$ws.Range("I117").Value = 748
$ws.Range("J117").Value = 499
$ws.Range("K117").Value = 2244
$ws.Range("L117").Value = 1497
$ws.Range("M117").Value = 1198
$ws.Range("N117").Value = -8381

$ws.Range("H129").Value = 3927
$ws.Range("J129").Value = 3927
$ws.Range("L129").Value = 11781
$ws.Range("N129").Value = -21781

$ws.Range("H136").Value = 890
$ws.Range("I136").Value = 890
$ws.Range("K136").Value = 2670
$ws.Range("M136").Value = 2430

$ws.Range("H141").Value = 250
$ws.Range("I141").Value = 250
$ws.Range("K141").Value = 750
$ws.Range("M141").Value = 4430

$ws = $wb.Sheets.Item("GSM")
$ws.Range("H97").Value = 1313.6364
$ws.Range("I97").Value = 1283.3334
$ws.Range("K97").Value = 1283.3334
$ws.Range("M97").Value = -787.3334

$ws.Range("H102").Value = 3803.6
$ws.Range("I102").Value = 3629.5
$ws.Range("J102").Value = 4500
$ws.Range("K102").Value = 3629.5
$ws.Range("L102").Value = 4500
$ws.Range("M102").Value = -2007.5
$ws.Range("N102").Value = -7744

$ws.Range("H114").Value = 97994.60000000001
$ws.Range("J114").Value = 97994.60000000001
$ws.Range("L114").Value = 97994.60000000001
$ws.Range("N114").Value = -106672.6

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H46").Value = 2706.7144
$ws.Range("I46").Value = 1200
$ws.Range("J46").Value = 2957.8333
$ws.Range("K46").Value = 1200
$ws.Range("L46").Value = 2957.8333
$ws.Range("M46").Value = -1012
$ws.Range("N46").Value = -3333.8333

$ws.Range("H136").Value = 3486.1667
$ws.Range("I136").Value = 3302.4
$ws.Range("J136").Value = 4405
$ws.Range("K136").Value = 9907.200000000001
$ws.Range("L136").Value = 13215
$ws.Range("M136").Value = -7357.200000000001
$ws.Range("N136").Value = -18315

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H75").Value = 18000
$ws.Range("I75").Value = 18000
$ws.Range("K75").Value = 18000
$ws.Range("M75").Value = -17064

$ws.Range("H78").Value = 18000
$ws.Range("I78").Value = 18000
$ws.Range("K78").Value = 54000
$ws.Range("M78").Value = -49320

$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

$ws.Range("H122").Value = 2859.6
$ws.Range("I122").Value = 2324.5
$ws.Range("K122").Value = 6973.5
$ws.Range("M122").Value = -4523.5

$ws.Range("H132").Value = 3634.3333
$ws.Range("I132").Value = 3634.3333
$ws.Range("K132").Value = 10902.9999
$ws.Range("M132").Value = -8372.999899999999

$ws.Range("H135").Value = 47000
$ws.Range("J135").Value = 56000
$ws.Range("L135").Value = 56000
$ws.Range("N135").Value = -66140
